# Apply the edit described by the diff:
#  1. Row 5 (CIPROFAR 500MG 10 F.C.TAB): update "current balance" (H) and
#     "transactions" (N) values.
#  2. Insert a new item row ("حنه جلوري  بني 1 كيس") right after row 24
#     (جنتيانا نقط) / before the former row 25 (سرنجات 5 سم), pushing the
#     remaining rows (سرنجات 5 سم, كريم فاتيكا 125 مل, the totals row and
#     the footer row) down by one.
#  3. Renumber the "م" (sequence) column for the rows that moved.
#  4. Update the running total and restore the row heights used by the
#     totals / footer rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Row 5 content fix ------------------------------------------------
$ws.Range("H5").Value = "0:0"
$ws.Range("N5").Value = "2:0"

# --- 2. Insert the new row at position 25 --------------------------------
$ws.Rows("25").Insert()

# Copy the formatting (styles, borders, fills) of the row above so the new
# row matches the rest of the table exactly, then fix up the row height
# (Insert() does not carry the custom height along).
$ws.Range("A24:N24").Copy()
$ws.Range("A25:N25").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false
$ws.Rows("25").RowHeight = 25.5

# Re-create the merged cells for the new row.
$ws.Range("B25:G25").Merge()
$ws.Range("H25:K25").Merge()
$ws.Range("L25:M25").Merge()

# Fill in the new row's values.
$ws.Range("A25").Value = 22
$ws.Range("B25").Value = "حنه جلوري  بني 1 كيس "
$ws.Range("H25").Value = "4:0"
$ws.Range("L25").Value = 40
$ws.Range("N25").Value = "1:0"

# --- 3. Renumber the rows that shifted down ------------------------------
$ws.Range("A26").Value = 23
$ws.Range("A27").Value = 24

# --- 4. Fix up the totals row and footer row -----------------------------
$ws.Rows("28").RowHeight = 25.5
$ws.Range("K28").Value = 1299.1400000000001

$ws.Rows("29").RowHeight = 17.25

Write-Host "edit complete"
